$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This worksheet is a weekly price log. A new weekly record is inserted at
# row 87 (pushing all subsequent records down by one row), and the record
# that used to be last (row 115) becomes the new last record at row 116.

$firstRow = 87
$lastRowOld = 115
$lastRowNew = 116
$lastCol = 18   # column R

# 1) Capture the original values of rows 87..115 (columns A..R) before we
#    start overwriting anything.
$data = @{}
for ($r = $firstRow; $r -le $lastRowOld; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals += ,$ws.Cells.Item($r, $c).Value2
    }
    $data[$r] = $rowVals
}

# 2) Shift every old row r (87..114) down into new row r+1 (88..115).
for ($r = $lastRowOld - 1; $r -ge $firstRow; $r--) {
    $srcVals = $data[$r]
    $destRow = $r + 1
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c - 1]
    }
}

# 3) Old row 115 becomes new row 116.
$srcVals = $data[$lastRowOld]
for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item($lastRowNew, $c).Value = $srcVals[$c - 1]
}

# 4) New row 87 keeps the same record as the old row 87 except for the date
#    (column D / index 4), which becomes 44468.
$newRow87 = $data[$firstRow]
for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item($firstRow, $c).Value = $newRow87[$c - 1]
}
$ws.Cells.Item($firstRow, 4).Value = 44468

# Keep the date column formatted the same way as the rest of column D.
$ws.Range("D116").NumberFormat = $ws.Range("D115").NumberFormat
